$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving exact formatting
# (no numeric/date coercion, no leading/trailing whitespace trimming),
# while leaving the cell's style untouched (matches original inlineStr
# cells which carry no explicit style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "51.626.89"
Set-TextValue $ws.Range("E2") "  -0.74%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.782.23"
Set-TextValue $ws.Range("E3") "  -0.38%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "352.63"
Set-TextValue $ws.Range("E5") "  -1.81%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "109.02"
Set-TextValue $ws.Range("E6") "  -0.76%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.550"
Set-TextValue $ws.Range("E7") "  -2.58%  "

# Row 8 - USDC
Set-TextValue $ws.Range("E8") "  +0.03%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("E9") "  +2.59%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "39.69"
Set-TextValue $ws.Range("E10") "  -1.07%  "

# Row 11 - TRON
Set-TextValue $ws.Range("E11") "  +2.47%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("E12") "  -2.32%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("E13") "  +3.06%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.67"
Set-TextValue $ws.Range("E14") "  +1.10%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.220.68"
Set-TextValue $ws.Range("E15") "  -0.31%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.803.47"
Set-TextValue $ws.Range("E16") "  +0.31%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.927"
Set-TextValue $ws.Range("E17") "  -1.76%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "51.615.04"
Set-TextValue $ws.Range("E18") "  -0.62%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "7.74"
Set-TextValue $ws.Range("E19") "  +4.35%  "

# Row 20 - ImmutableX
Set-TextValue $ws.Range("D20") "3.16"
Set-TextValue $ws.Range("E20") "  +0.58%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "13.16"
Set-TextValue $ws.Range("E21") "  +0.94%  "

# Row 22 - ShibaInu
Set-TextValue $ws.Range("D22") "0.0₃0963"
Set-TextValue $ws.Range("E22") "  -2.12%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("E23") "  -0.57%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "267.41"
Set-TextValue $ws.Range("E24") "  -2.39%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.73"
Set-TextValue $ws.Range("E25") "  -0.92%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("D26") "26.12"
Set-TextValue $ws.Range("E26") "  -2.07%  "

# Row 27 - Dai
Set-TextValue $ws.Range("D27") "0.998"
Set-TextValue $ws.Range("E27") "  -0.15%  "

# Row 28 - Kaspa
Set-TextValue $ws.Range("E28") "  +13.25%  "

# Row 29 - Cosmos
Set-TextValue $ws.Range("E29") "  +0.50%  "

# Row 30 - InjectiveProtocol
Set-TextValue $ws.Range("D30") "37.09"
Set-TextValue $ws.Range("E30") "  +7.50%  "

# Row 31 - Toncoin
Set-TextValue $ws.Range("E31") "  -1.35%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("E32") "  +7.11%  "

# Row 33 - OKB
Set-TextValue $ws.Range("D33") "51.68"
Set-TextValue $ws.Range("E33") "  +0.00%  "

# Row 34 - RenderToken
Set-TextValue $ws.Range("D34") "5.66"
Set-TextValue $ws.Range("E34") "  +7.58%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("D36") "0.0832"
Set-TextValue $ws.Range("E36") "  -1.82%  "

# Row 37 - FirstDigitalUSD
Set-TextValue $ws.Range("E37") "  -0.03%  "

# Row 38 - Celestia
Set-TextValue $ws.Range("D38") "18.47"
Set-TextValue $ws.Range("E38") "  +1.34%  "

# Row 39 - LidoDAOToken
Set-TextValue $ws.Range("E39") "  -2.97%  "

# Row 40 - ARBITRUM
Set-TextValue $ws.Range("E40") "  -1.69%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("E41") "  -1.49%  "

# Row 42 - Stellar
Set-TextValue $ws.Range("E42") "  -0.78%  "

# Row 43 - was Monero, becomes EnergySwap (rows 43/44 swap coin identity)
Set-TextValue $ws.Range("B43") "EnergySwap"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D43") "22.18"
Set-TextValue $ws.Range("E43") "  -0.58%  "

# Row 44 - was EnergySwap, becomes Monero
Set-TextValue $ws.Range("B44") "Monero"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D44") "120.27"
Set-TextValue $ws.Range("E44") "  -1.90%  "

# Row 45 - WEMIXToken
Set-TextValue $ws.Range("D45") "2.17"
Set-TextValue $ws.Range("E45") "  -3.50%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "2.126.90"
Set-TextValue $ws.Range("E46") "  +1.79%  "

# Row 47 - NEARProtocol
Set-TextValue $ws.Range("E47") "  +1.80%  "

# Row 48 - ApeXProtocol
Set-TextValue $ws.Range("E48") "  +4.31%  "

# Row 49 - THORChain
Set-TextValue $ws.Range("D49") "5.42"
Set-TextValue $ws.Range("E49") "  -5.14%  "

# Row 50 - SEI
Set-TextValue $ws.Range("D50") "0.904"
Set-TextValue $ws.Range("E50") "  -3.00%  "

# Row 51 - TrustWalletToken
Set-TextValue $ws.Range("E51") "  +8.57%  "
